# Fgf18-Fgfr4.xlsx was regenerated from an updated (TPM-based) NATMI run.
# The "ECs" sending-cluster block (old rows 2-4) is no longer part of the
# output, so those three rows are removed and the remaining "FAPs"/"MuSCs"
# sending-cluster rows shift up to become the new rows 2-7. All of the
# downstream, TPM-derived statistics (receptor stats + edge weights/
# specificities, columns I-T) are refreshed with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three "ECs" sending-cluster rows; everything below shifts up.
$ws.Range("A2:A4").EntireRow.Delete()

# New values for columns I-T (Receptor stats through Edge specificities)
# for each of the remaining rows, now renumbered 2-7.
$updates = @{
    2  = @{ I = 0.9145494540267;     J = 0.9145494540267;     K = 2; L = 0.6666666666666666;
            M = 0.08378199999999998; N = 0.251346;             O = 0.007571394704126512;  P = 0.007571394704126512;
            Q = 0.7348381258973331;  R = 6.613543133075999;    S = 0.006924414892879549;  T = 0.006924414892879549 }
    3  = @{ I = 0.9145494540267;     J = 0.9145494540267;     K = 1; L = 0.3333333333333333;
            M = 0.07352966666666667; N = 0.220589;             O = 0.006644889460697858;  P = 0.006644889460697857;
            Q = 0.6449165984482222;  R = 5.804249386034;       S = 0.006077080028348999;  T = 0.006077080028348998 }
    4  = @{ I = 0.9145494540267;     J = 0.9145494540267;     K = 3; L = 1;
            M = 10.908285;           N = 32.724855;            O = 0.9857837158351757;    P = 0.9857837158351755;
            Q = 95.67477150406998;   R = 861.07294353663;      S = 0.9015479591054715;    T = 0.9015479591054714 }
    5  = @{ I = 0.08545054597330007; J = 0.08545054597330005; K = 2; L = 0.6666666666666666;
            M = 0.08378199999999998; N = 0.251346;             O = 0.007571394704126512;  P = 0.007571394704126512;
            Q = 0.06865929314533331; R = 0.6179336383079999;   S = 0.0006469798112469632; T = 0.000646979811246963 }
    6  = @{ I = 0.08545054597330007; J = 0.08545054597330005; K = 1; L = 0.3333333333333333;
            M = 0.07352966666666667; N = 0.220589;             O = 0.006644889460697858;  P = 0.006644889460697857;
            Q = 0.06025751281355556; R = 0.542317615322;       S = 0.0005678094323488594; T = 0.0005678094323488593 }
    7  = @{ I = 0.08545054597330007; J = 0.08545054597330005; K = 3; L = 1;
            M = 10.908285;           N = 32.724855;            O = 0.9857837158351757;    P = 0.9857837158351755;
            Q = 8.939332285309998;   R = 80.45399056778999;    S = 0.08423575672970425;   T = 0.08423575672970422 }
}

foreach ($rowNum in $updates.Keys) {
    $rowValues = $updates[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}
